# ---------------------------------------------------------------------------
# FreeCRMData.xlsx edit: "This is my second commit"
#
# Summary of the change (reconstructed from the OOXML diff):
#   * CustomerDetails: three customer rows get new names / mobiles / emails
#     (Ramakrishna/Suhas/Venkat -> Ramki/Suri/Jhon, new mobile numbers, new
#     email addresses).
#   * The NewAccount / DeleteAcc / Deposit sheets are swapped out for fresh
#     data and renamed: NewAccount -> "Deposit" (new deposit rows),
#     DeleteAcc -> "NewAcc1" (new account rows), Deposit -> "DelAcc" (new
#     delete-account row).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) CustomerDetails: update the three customer records.
# ---------------------------------------------------------------------------
$cust = $wb.Worksheets.Item("CustomerDetails")

$cust.Range("A2").Value = "Ramki"
$cust.Range("A3").Value = "Suri"
$cust.Range("A4").Value = "Jhon"

$cust.Range("H2").Value = "ramki1@gmail.com"
$cust.Range("H3").Value = "suri12@gmail.com"
$cust.Range("H4").Value = "Jhon@gmail.com"

$cust.Range("G2").Value = 64648
$cust.Range("G3").Value = 64649
$cust.Range("G4").Value = 64610

# Selection moves down one row (H4 -> H5) - do this before we touch the
# other sheets' selections so the very last Select() in the script is the
# one that determines the workbook's final ActiveSheet/ActiveTab.
$cust.Range("H5").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2) Rebuild the NewAccount / DeleteAcc / Deposit sheets with fresh data,
#    copying the existing header / data cell formatting so no new style
#    entries are introduced.
# ---------------------------------------------------------------------------
$wsNewAccount = $wb.Worksheets.Item("NewAccount")
$wsDeleteAcc  = $wb.Worksheets.Item("DeleteAcc")
$wsDeposit    = $wb.Worksheets.Item("Deposit")
$wsLogin      = $wb.Worksheets.Item("GuruLoginDetails1")

# --- header / data format sources (unchanged styles already in the file) ---
# Deposit!A1 carries the shaded header style used by all three grids.
# GuruLoginDetails1!B2 carries the centered, bordered data-row style.

# 2a) NewAccount sheet -> becomes "Deposit" (AccNo / Amount / Desc grid)
$wsNewAccount.Cells.Clear() | Out-Null
$wsDeposit.Range("A1").Copy() | Out-Null
$wsNewAccount.Range("A1:C1").PasteSpecial(-4122) | Out-Null
$wsLogin.Range("B2").Copy() | Out-Null
$wsNewAccount.Range("A2:C3").PasteSpecial(-4122) | Out-Null

$wsNewAccount.Range("A1").Value = "AccNo"
$wsNewAccount.Range("B1").Value = "Amount"
$wsNewAccount.Range("C1").Value = "Desc"
$wsNewAccount.Range("A2").Value = 89458
$wsNewAccount.Range("B2").Value = 1000
$wsNewAccount.Range("C2").Value = "BillPay"
$wsNewAccount.Range("A3").Value = 89459
$wsNewAccount.Range("B3").Value = 1000
$wsNewAccount.Range("C3").Value = "BillPay"

$wsNewAccount.Range("A2").Select() | Out-Null

# 2b) DeleteAcc sheet -> becomes "NewAcc1" (CustID / Deposite grid)
$wsDeleteAcc.Cells.Clear() | Out-Null
$wsDeposit.Range("A1").Copy() | Out-Null
$wsDeleteAcc.Range("A1:B1").PasteSpecial(-4122) | Out-Null
$wsLogin.Range("B2").Copy() | Out-Null
$wsDeleteAcc.Range("A2:B3").PasteSpecial(-4122) | Out-Null

$wsDeleteAcc.Range("A1").Value = "CustID"
$wsDeleteAcc.Range("B1").Value = "Deposite"
$wsDeleteAcc.Range("A2").Value = 34189
$wsDeleteAcc.Range("B2").Value = 5000
$wsDeleteAcc.Range("A3").Value = 62819
$wsDeleteAcc.Range("B3").Value = 6000

# 2c) Deposit sheet -> becomes "DelAcc" (AccountNo grid)
$wsDeposit.Cells.Clear() | Out-Null
$wsDeposit.Range("A1").Value = "AccountNo"
$wsDeposit.Range("A2").Value = 89666

$wsLogin.Range("B2").Copy() | Out-Null
$wsDeposit.Range("A2").PasteSpecial(-4122) | Out-Null
$wsDeposit.Range("A1").Value = "AccountNo"
$wsDeposit.Range("A2").Value = 89666

# Re-apply the shaded header style on the now-single header cell.
$wsNewAccount.Range("A1").Copy() | Out-Null
$wsDeposit.Range("A1").PasteSpecial(-4122) | Out-Null
$wsDeposit.Range("A1").Value = "AccountNo"

$wsDeposit.Range("A1:A2").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3) Rename the three sheets into their final names (via a temp name to
#    avoid collisions while the rename is in flight).
# ---------------------------------------------------------------------------
$wsDeposit.Name = "TempDeposit"
$wsNewAccount.Name = "Deposit"
$wsDeleteAcc.Name = "NewAcc1"
$wsTempDeposit = $wb.Worksheets.Item("TempDeposit")
$wsTempDeposit.Name = "DelAcc"

# ---------------------------------------------------------------------------
# 4) Final selection: NewAcc1 (old DeleteAcc) is the active tab in the
#    source workbook (activeTab="4"), selected at D8. This must be the very
#    last Select() call so it "wins" as the workbook's active sheet.
# ---------------------------------------------------------------------------
$wsNewAcc1 = $wb.Worksheets.Item("NewAcc1")
$wsNewAcc1.Range("D8").Select() | Out-Null
